$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "anisa@gmail.com"
$ws.Range("B2").Value = "xyz123"

$ws.Range("C3").Select()
